$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the summary block ---
# Total "VALOR MORA" (E11) grew because more records were added.
$ws.Range("E11").Value = 1180820
# "Cant. Trabajadores" (C13) and "Cant. Periodos" (F13) counts increased.
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 3

# --- Make room for the new data rows ---
# The table originally had 8 data rows (16-23); the new data has 13 rows
# (16-28). Insert 5 blank rows above the old last row (23) so that row 23
# (with its special "closing" bottom-border style) slides down intact to
# become the new last row (28), and the footer block shifts from rows
# 28-29 down to rows 33-34.
$ws.Rows("23:27").Insert()

# Give the newly inserted rows (23-27) the same formatting as a normal
# interior data row (row 22) rather than the default blank formatting.
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J27").PasteSpecial(-4122)

# --- Fill in the new worker account-statement data (rows 16-28) ---
$data = @(
    @(16, "CC", "1043963669", "SERGIO ADRIAN ORTEGA CARABALLO", "2507", 100000, 2500000),
    @(17, "CC", "1043963669", "SERGIO ADRIAN ORTEGA CARABALLO", "2506", 100000, 2500000),
    @(18, "CC", "1043963669", "SERGIO ADRIAN ORTEGA CARABALLO", "2505", 100000, 2500000),
    @(19, "CC", "1051675065", "DAVID ERNESTO VIDES JIMENEZ", "2507", 160000, 4000000),
    @(20, "CC", "1032482224", "LAURA ALEJANDRA ESPITIA GOMEZ", "2506", 160000, 4000000),
    @(21, "CC", "1032482224", "LAURA ALEJANDRA ESPITIA GOMEZ", "2505", 160000, 4000000),
    @(22, "CC", "1051675065", "DAVID ERNESTO VIDES JIMENEZ", "2507", 14000, 3500000),
    @(23, "CC", "1050967671", "JOEMIS PATRICIA ARNEDO GAMBIN", "2507", 56940, 1423500),
    @(24, "CC", "1050967671", "JOEMIS PATRICIA ARNEDO GAMBIN", "2506", 56940, 1423500),
    @(25, "CC", "1050967671", "JOEMIS PATRICIA ARNEDO GAMBIN", "2505", 56940, 1423500),
    @(26, "CC", "1007848449", "CHELSEA STELLA BERMUDEZ RUEDA", "2507", 72000, 1800000),
    @(27, "CC", "1007848449", "CHELSEA STELLA BERMUDEZ RUEDA", "2506", 72000, 1800000),
    @(28, "CC", "1007848449", "CHELSEA STELLA BERMUDEZ RUEDA", "2505", 72000, 1800000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
